# Updated cryptos list on Mon Apr  3 21:11:25 UTC 2023 with GitHub Actions
# Re-applies the scraped coinranking.com snapshot: refreshed Price/Volume(1h)
# figures for every row, plus three ranking swaps (Litecoin/ShibaInu,
# Filecoin/HuobiToken, Decentraland/PancakeSwap) that changed places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.616.44"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "1.791.17"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "'1.006"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'304.86"
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("D7").Value = "'0.4920"
$ws.Range("E7").Value = "  -5.44%  "
$ws.Range("D8").Value = "'0.3811"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.09137"
$ws.Range("E9").Value = "  +14.26%  "
$ws.Range("D10").Value = "'1.085"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "'40.51"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").Value = "'1.008"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "'6.206"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "'20.30"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "1.797.06"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "'7.117"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001099"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'91.52"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'0.06559"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'16.93"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").Value = "'5.900"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "27.679.94"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "'10.90"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").Value = "'2.217"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").Value = "'157.38"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").Value = "2.003.87"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "'20.24"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "'2.360"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "'126.36"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("D31").Value = "'0.1064"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "'1.040"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.620"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.470"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "'0.06749"
$ws.Range("E35").Value = "  -6.87%  "
$ws.Range("D36").Value = "'8.747"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "'0.02277"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "'0.2111"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").Value = "'11.22"
$ws.Range("E39").Value = "  -7.03%  "
$ws.Range("D40").Value = "'4.876"
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("D41").Value = "'0.6058"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").Value = "'1.006"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "'1.136"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").Value = "'12.85"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.646"
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5763"
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D48").Value = "'123.26"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "'1.908"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").Value = "'1.158"
$ws.Range("E50").Value = "  -5.17%  "
$ws.Range("D51").Value = "'0.06699"
$ws.Range("E51").Value = "  -0.70%  "
